# Fix all ICDC breed testcases
# Applies the TC45 Canine Filter Breed-Weimrnr worksheet update:
#  - refresh the CasesTab / SamplesTab / FilesTab Neo4j queries
#  - add a new StudyFilesTab row
#  - adjust row heights / selection to match the refreshed content

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (CasesTab): new case query, with apoc age/weight integer coalesce ---
$ws.Range("B2").Value = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
 WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
WHERE demo.breed IN [''Weimaraner'']
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
  coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '''') AS Age,
       coalesce(demo.sex, '''') AS Sex,
       coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '''') AS `Weight (kg)`,
       coalesce(diag.best_response, '''') AS `Response to Treatment`,
       coalesce(co.cohort_description, '''') AS `Cohort`
order by c.case_id asc
limit 100'

# --- Row 3 (SamplesTab): new sample query, with explicit ORDER BY / LIMIT ---
$ws.Range("B3").Value = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE demo.breed IN [''Weimaraner''] 
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '''') AS `Sample ID`, 
        coalesce(c.case_id, '''') AS `Case ID`, 
        coalesce(demo.breed,'''') AS Breed,
        coalesce(diag.disease_term,'''') AS Diagnosis, 
        coalesce(samp.sample_site, '''') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '''') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '''') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '''') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '''') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '''') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '''') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '''') AS `Sample Preservation`
Order by samp.sample_id LIMIT 100'

# --- Row 4 (FilesTab): new case-files query, with size unit computation ---
$ws.Range("B4").Value = 'MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f)-[*]->(samp:sample)
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN [''Weimaraner''] 
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN 
        coalesce(f.file_name, '''') AS `File Name`,
        coalesce(f.file_format, '''') AS `Format`,
        coalesce(f.file_type, '''') AS `File Type`,
       CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
        coalesce(labels(parent)[0], '''') AS `Association`,
        coalesce(f.file_description, '''') AS `Description`,
        coalesce(samp.sample_id, '''') AS `Sample ID`,
        coalesce(c.case_id, '''') AS `Case ID`,
        coalesce(demo.breed,'''') AS Breed ,
        coalesce(diag.disease_term,'''') AS Diagnosis
Order By f.file_name LIMIT 100'

# --- Row 5 (StudyFilesTab): brand-new row for the study-files query ---
$ws.Range("A5").Value = 'StudyFilesTab'
$ws.Range("B5").Value = 'MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
WHERE demo.breed IN [''Weimaraner''] 
WITH DISTINCT f,  s, c, demo, diag
WITH
        f, c, demo, diag, s,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH    
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH    
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '''') AS `File Name`,
  coalesce(f.file_type, '''') AS `File Type`,
  coalesce("study", '''') AS `Association`,
  coalesce(f.file_description, '''') AS `Description`,
  coalesce(f.file_format, '''') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
  coalesce(s.clinical_study_designation,'''') AS `Study Code`
Order By f.file_name LIMIT 100'
$ws.Range("C5").Value = $ws.Range("C4").Value()
$ws.Range("D5").Value = $ws.Range("D4").Value()
$ws.Range("E5").Value = $ws.Range("E4").Value()

$ws.Range("B5").WrapText = $true
$ws.Range("C5").WrapText = $true

# --- Row heights to match the refreshed (longer) query text ---
$ws.Rows.Item(2).RowHeight = 288
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 409.6
$ws.Rows.Item(5).RowHeight = 374.4

# --- Selection / view state ---
$ws.Range("C4:E5").Select()
